$p = $ppt.ActivePresentation

# --- 1. Slide 6: change the table's style (Table Design gallery pick) ---
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{56CC2776-F674-44E6-9EAF-9412ACCE41F9}")

# --- 2. Switch the presentation's theme colors from "Integral" to the
#        default "Office Theme" palette (Design tab -> Themes gallery).
#        ThemeColorScheme index order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
#        5-10 accent1-6, 11 hlink, 12 folHlink.
$master = $p.Designs.Item(1).SlideMaster
$colors = $master.Theme.ThemeColorScheme

$colors.Item(3).RGB  = 0x6A5444   # dk2      44546A
$colors.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$colors.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$colors.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$colors.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$colors.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$colors.Item(9).RGB  = 0xC47244   # accent5  4472C4
$colors.Item(10).RGB = 0x47AD70   # accent6  70AD47
$colors.Item(11).RGB = 0xC16305   # hlink    0563C1
$colors.Item(12).RGB = 0x724F95   # folHlink 954F72
